$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in test case data (rows 6-10, columns B-I) ---
$ws.Range("B6").Value = 8
$ws.Range("C6").Value = 11
$ws.Range("D6").Value = 27
$ws.Range("E6").Value = 337131148
$ws.Range("F6").Value = 5
$ws.Range("G6").Formula = "=INT((1/B6-1/C6+1/D6)*(31536000)*F6)"
$ws.Range("H6").Formula = "=G6+E6"
$ws.Range("I6").Formula = '=IF(G6>0, "Increase", "Decrease")'

$ws.Range("B7").Value = 3
$ws.Range("C7").Value = 11
$ws.Range("D7").Value = 27
$ws.Range("E7").Value = 337131148
$ws.Range("F7").Value = 5
$ws.Range("G7").Formula = "=INT((1/B7-1/C7+1/D7)*(31536000)*F7)"
$ws.Range("H7").Formula = "=G7+E7"
$ws.Range("I7").Formula = '=IF(G7>0, "Increase", "Decrease")'

$ws.Range("B8").Value = 8
$ws.Range("C8").Value = 11
$ws.Range("D8").Value = 3
$ws.Range("E8").Value = 337131148
$ws.Range("F8").Value = 5
$ws.Range("G8").Formula = "=INT((1/B8-1/C8+1/D8)*(31536000)*F8)"
$ws.Range("H8").Formula = "=G8+E8"
$ws.Range("I8").Formula = '=IF(G8>0, "Increase", "Decrease")'

$ws.Range("B9").Value = 8
$ws.Range("C9").Value = 1
$ws.Range("D9").Value = 27
$ws.Range("E9").Value = 337131148
$ws.Range("F9").Value = 5
$ws.Range("G9").Formula = "=INT((1/B9-1/C9+1/D9)*(31536000)*F9)"
$ws.Range("H9").Formula = "=G9+E9"
$ws.Range("I9").Formula = '=IF(G9>0, "Increase", "Decrease")'

$ws.Range("B10").Value = 30
$ws.Range("C10").Value = 11
$ws.Range("D10").Value = 100
$ws.Range("E10").Value = 337131148
$ws.Range("F10").Value = 5
$ws.Range("G10").Formula = "=INT((1/B10-1/C10+1/D10)*(31536000)*F10)"
$ws.Range("H10").Formula = "=G10+E10"
$ws.Range("I10").Formula = '=IF(G10>0, "Increase", "Decrease")'
